# Feb 13: Adding line chart, data manipulation
# (Per the source diff, the persisted change is the new data row for Feb 13;
#  add it with the same date/number formatting used by the preceding rows.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entry: 2022-02-13 (serial 44605), 1.5 hours, progress note.
$ws.Range("A8").Value = 44605
$ws.Range("A8").NumberFormat = $ws.Range("A7").NumberFormat

$ws.Range("B8").Value = 1.5
$ws.Range("B8").HorizontalAlignment = $ws.Range("B7").HorizontalAlignment
$ws.Range("B8").VerticalAlignment = $ws.Range("B7").VerticalAlignment

$ws.Range("C8").Value = "Adding a line chart to help understand trends over time"

# Match the active selection left behind by the edit.
$ws.Range("C8").Select()
